$wb = $excel.ActiveWorkbook

# Rename "Paineis DARQ" to "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" to "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete "Desarquivamentos Pendentes" sheet
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true

# Restore original active sheet/tab selection
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
